# Add IBAN values (column M) for each employee data row (rows 2-104) on "Hoja1".
# Column M already has its header ("IBAN", row 1) in place; this fills in the per-row
# account numbers, which introduces the corresponding new shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ibanValues = @("ES3520960043043554600000","ES1132145464978452163421","DK1531645124013461205164","ES0465614874115615445616","RO2732569523636220165156","DE5324561937321546497521","MC8636520125478451012515","GR5020125003625201112544","ES6821651484840980008984","FI5620960043043554600000","ES7321564975463245467995","LT2932566221452587754554","EE8723215465445456411515","BE5900750184300702510000","ES8421584976342154655487","SM5725894363125485700145","ES4296431245768150005156","AT6125030000694574745458","IT4415953684521254695203","ES5520960043013096200000","DK2300750184300702510000","ES2223455254353263234457","GR4920910936583000000000","ES0220960043022159000000","DE9712669681515112121210","FI4812121455156115515151","ES6956187775465550000651","ES0225516848741156151054","PT2164578946120051516490","ES5234698752134600549403","ES0766649444572310000255","FR9523185484005641685100","ES1152198484632100515144","CZ5836250012544785523365","AT6722515651565640081000","DE5421508149675421346497","DE8421346154203164978451","ES9825187786651225455548","ES1223164897632213030615","ES1396536214585214585214","ES2785461325461978750005","FI7724587946012003165464","ES8520960043083071400000","ES6320960043092158800000","ES1821654587165156484454","ES6751651681001210656510","ES3866552211448855332200","GB9720910936583000000000","DE5801821135900205540000","DE6822631245246916432102","ES1320960043013075700000","SM9425635478521002541225","ES5032154697855423121000","GR5836521452326500658485","GB0420008521458775113366","ES0751651487560005118185","ES7320036525587451152165","ES9001826530120201560000","ES6421651651472511133551","ES1020960043013000100000","GB3736585214520025478551","ES2212548523695214585214","ES4431624561312546920007","ES6036154231642500312566","ES7744875664957231645789","ES7920960031442124800000","ES0933620012547852100256","ES7533218885651445121022","ES2462581542003690044508","ES7325165151488666365100","ES1020960043013000100000","PT8336952365200014425254","ES3365168874651561561500","ES4020960583861234500000","ES2621416325881510005514","LU0232628484114115151115","ES7732562000125478520015","ES3020960056163231500000","ES5463516541848944000984","ES4623658965214585223202","FI1232658012547712548745","ES9623652365232254222000","FR5720012541100023365233","ES7732584216551684051000","ES2855065688611051056105","ES6095485212145484010000","LT6321856333256985542360","ES4836245978463245679001","ES8931245164326597845124","SM5723221158422545471411","SE5232574512025411002255","ES6020960043073468900000","ES2331215643125060225021","AT7485550564156165145610","ES4165165654648886005001","ES4926221011898048788896","ES8112548521458742146695","IE0951556584121251000254","DK1832541112541220000588","LT7862541122001110105611","AT8065645150005168448896","IT8626551681877651415636","HU2299558741226555551120")

$startRow = 2
for ($i = 0; $i -lt $ibanValues.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Range("M" + $row)
    $cell.Value = $ibanValues[$i]
    # New cells pick up no explicit formatting of their own (only the column's
    # default style applies), matching the rest of the sheet's data cells.
    $cell.Style = "Normal"
}
